# Update "想去人数" (interest count, column F) and "最低票价" (min ticket
# price, column G) figures for several 漫展 (con) listings, mirrored across
# the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("G2").Value  = 55    # 南昌·云汐动漫游戏嘉年华 (min ticket price 45 -> 55)
$ws1.Range("G3").Value  = 65    # 抚州·BM次元盛典运动番only (min ticket price 55 -> 65)
$ws1.Range("F5").Value  = 1841  # 南昌·ACG CLUB动漫游戏嘉年华 (interest 1840 -> 1841)
$ws1.Range("F9").Value  = 2459  # 南昌·CM02动漫游戏博览会 (interest 2456 -> 2459)
$ws1.Range("F11").Value = 80    # 九江·第四届ACD动漫游戏嘉年华 (interest 79 -> 80)
$ws1.Range("F15").Value = 38    # 宜春·静卿缤纷仲夏国风动漫文化展 (interest 37 -> 38)
$ws1.Range("F18").Value = 20    # 上饶·宅舞联萌·随舞动漫派对 (interest 19 -> 20)
$ws1.Range("F26").Value = 1540  # 江西·次元星河国风动漫游戏嘉年华 (interest 1533 -> 1540)
$ws1.Range("F27").Value = 20    # 赣州·明日方舟only叙拉古夜宴3.0 (interest 19 -> 20)
$ws1.Range("F29").Value = 385   # 景德镇·第十五届瓷都ACG动漫游戏博览会 (interest 379 -> 385)
$ws1.Range("F30").Value = 193   # 九江·第一届异次元动漫嘉年华 (interest 190 -> 193)
$ws1.Range("F31").Value = 291   # 南昌·第一届异次元动漫嘉年华 (interest 290 -> 291)
$ws1.Range("F32").Value = 391   # 赣州·第二届异次元动漫嘉年华 (interest 388 -> 391)

# --- Sheet "全部类型" (same listings mirrored, with a one-off row shift) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("G2").Value  = 55    # 南昌·云汐动漫游戏嘉年华 (min ticket price 45 -> 55)
$ws4.Range("G3").Value  = 65    # 抚州·BM次元盛典运动番only (min ticket price 55 -> 65)
$ws4.Range("F5").Value  = 1841  # 南昌·ACG CLUB动漫游戏嘉年华 (interest 1840 -> 1841)
$ws4.Range("F10").Value = 2459  # 南昌·CM02动漫游戏博览会 (interest 2456 -> 2459)
$ws4.Range("F12").Value = 80    # 九江·第四届ACD动漫游戏嘉年华 (interest 79 -> 80)
$ws4.Range("F16").Value = 38    # 宜春·静卿缤纷仲夏国风动漫文化展 (interest 37 -> 38)
$ws4.Range("F19").Value = 20    # 上饶·宅舞联萌·随舞动漫派对 (interest 19 -> 20)
$ws4.Range("F27").Value = 1540  # 江西·次元星河国风动漫游戏嘉年华 (interest 1534 -> 1540)
$ws4.Range("F28").Value = 20    # 赣州·明日方舟only叙拉古夜宴3.0 (interest 19 -> 20)
$ws4.Range("F30").Value = 385   # 景德镇·第十五届瓷都ACG动漫游戏博览会 (interest 379 -> 385)
$ws4.Range("F31").Value = 193   # 九江·第一届异次元动漫嘉年华 (interest 190 -> 193)
$ws4.Range("F32").Value = 291   # 南昌·第一届异次元动漫嘉年华 (interest 290 -> 291)
$ws4.Range("F33").Value = 391   # 赣州·第二届异次元动漫嘉年华 (interest 388 -> 391)

Write-Output "Updated 13 cells on sheet 1 and 13 cells on sheet 4."
